# Daily attendance processing - rotate the "Recorded By" contributor list
# in column G so the first-listed recorder moves to the end of the list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $text = $cell.Text

    if ($text -and $text.Contains(",")) {
        $parts = $text -split ","
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }
        $rotated = $trimmed[1..($trimmed.Length - 1)] + $trimmed[0]
        $cell.Value = [string]::Join(", ", $rotated)
    }
}
